# Auto-generated Excel COM-interop edit script
# Applies the diff: updates Situacao/Observacao for rows 583-584
# and appends 10 new sales rows (585-594) to Planilha2 / Tabela3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- Update existing rows 583 and 584: Situacao -> Instalada, Observacao updated ---

$ws.Cells.Item(583, 9).Value = "Instalada"
$j583 = @'
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 06/03/2025 13:58
Remarcado de 06/03/2025 de Tarde Para 07/03/2025 de Tarde
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 10/03/2025 17:41
Sua Venda mudou de Aprovisionamento Para Instalada

'@
Set-TextCell $ws.Cells.Item(583, 10) $j583
$ws.Rows.Item(583).EntireRow.AutoFit()

$ws.Cells.Item(584, 9).Value = "Instalada"
$j584 = @'
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 07/03/2025 11:49
Sua Venda mudou de Aprovisionamento Para Instalada

'@
Set-TextCell $ws.Cells.Item(584, 10) $j584
$ws.Rows.Item(584).EntireRow.AutoFit()

# --- Append new sales rows 585-594 ---

# Row 585
$r = 585
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '00826848125'
$ws.Cells.Item($r, 2).Value = 'Alzirene Cândida da Silva '
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '63992699286'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '06/03/2025 16:39'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '07/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Instalada'
$ws.Cells.Item($r, 10).Value = @'
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 07/03/2025 11:47
Sua Venda mudou de Aprovisionamento Para Instalada

'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 586
$r = 586
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '03993926617'
$ws.Cells.Item($r, 2).Value = 'ADRIANA FERREIRA DE JESUS ALVES'
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '31971009748'
$ws.Cells.Item($r, 4).Value = '400MB'
$ws.Cells.Item($r, 5).Value = 30
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '06/03/2025 18:48'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '07/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Instalada'
$ws.Cells.Item($r, 10).Value = @'
VERO
Contrato: 2344249   Venda #689064
Cód. endereço: 1487355
Inventário: 58377018
'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'SKY'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 587
$r = 587
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '03167641673'
$ws.Cells.Item($r, 2).Value = 'WESLLEY COSTA VILACA'
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '31982954636'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '06/03/2025 19:14'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '07/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Instalada'
$ws.Cells.Item($r, 10).Value = @'
VERO
Contrato: 2344249   Venda #689064
Cód. endereço: 1487355
Inventário: 58377018
'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'SKY'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 588
$r = 588
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '05284784278'
$ws.Cells.Item($r, 2).Value = 'ALICE MORAES GUIMARAES'
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '92993150726'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '10/03/2025 17:40'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '14/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Instalada'
$ws.Cells.Item($r, 10).Value = @'
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 14/03/2025 12:59
Remarcado de 13/03/2025 de Manhã Para 14/03/2025 de Manhã
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 17/03/2025 11:27
Sua Venda mudou de Aprovisionamento Para Instalada

'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 589
$r = 589
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '08197351759'
$ws.Cells.Item($r, 2).Value = 'Helena Rúbia Sena'
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '21965019039'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '13/03/2025 17:12'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '15/03/2025'
$ws.Cells.Item($r, 8).Value = 'Tarde'
$ws.Cells.Item($r, 9).Value = 'Instalada'
$ws.Cells.Item($r, 10).Value = @'
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 17/03/2025 11:25
Sua Venda mudou de Aprovisionamento Para Instalada

'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 590
$r = 590
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '72824247215'
$ws.Cells.Item($r, 2).Value = 'SILVIA GERONIMO GUEDES'
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '92981658506'
$ws.Cells.Item($r, 4).Value = '500MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '13/03/2025 18:13'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '17/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Aprovisionamento'
$ws.Cells.Item($r, 10).Value = ''
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 591
$r = 591
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '02052239316  '
$ws.Cells.Item($r, 2).Value = 'Alessandro da Silva  '
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '61993980524'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '14/03/2025 12:03'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '17/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Instalada'
$ws.Cells.Item($r, 10).Value = @'
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 14/03/2025 12:57
Remarcado de 15/03/2025 de Manhã Para 17/03/2025 de Manhã
--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 17/03/2025 11:27
Sua Venda mudou de Aprovisionamento Para Instalada

'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 592
$r = 592
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '03626363708'
$ws.Cells.Item($r, 2).Value = 'Adriana Da Silva Chaves Martins  '
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '24999324822'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '14/03/2025 14:39'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '15/03/2025'
$ws.Cells.Item($r, 8).Value = 'Tarde'
$ws.Cells.Item($r, 9).Value = 'Aprovisionamento'
$ws.Cells.Item($r, 10).Value = ''
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 593
$r = 593
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '60899409687'
$ws.Cells.Item($r, 2).Value = 'CLAUDIA MAGDA DA SILVA'
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '31988569121'
$ws.Cells.Item($r, 4).Value = '700MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '14/03/2025 15:46'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '15/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Aprovisionamento'
$ws.Cells.Item($r, 10).Value = @'
Dados da integração:
Contrato: 2360935
Cód. endereço: 24786877
Inventário: 17793439
Cód. agendamento: 10d6b121-53f7-4535-be51-96aa21f38646
Cód. documento: 20250314154003873
Produto: BL_780MB
Agendamento: 17/03/2025 08:00:00 à 17/03/2025 12:00:00
'@
$ws.Cells.Item($r, 11).Value = 'Base'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'SKY'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

# Row 594
$r = 594
$rowRange = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 14))
$rowRange.Style = "Normal"
$ws.Cells.Item($r, 1).NumberFormat = "@"
$ws.Cells.Item($r, 1).Value = '35943378120'
$ws.Cells.Item($r, 2).Value = 'Genoveva Angela pinto  '
$ws.Cells.Item($r, 3).NumberFormat = "@"
$ws.Cells.Item($r, 3).Value = '61993169189'
$ws.Cells.Item($r, 4).Value = '500MB'
$ws.Cells.Item($r, 5).Value = 35
$ws.Cells.Item($r, 6).NumberFormat = "@"
$ws.Cells.Item($r, 6).Value = '14/03/2025 17:31'
$ws.Cells.Item($r, 7).NumberFormat = "@"
$ws.Cells.Item($r, 7).Value = '16/03/2025'
$ws.Cells.Item($r, 8).Value = 'Manhã'
$ws.Cells.Item($r, 9).Value = 'Aprovisionamento'
$ws.Cells.Item($r, 10).Value = @'
sem agendamento--%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&----%$&--
     Atualizado em 14/03/2025 17:48
Remarcado de 15/03/2025 de Manhã Para 16/03/2025 de Manhã

'@
$ws.Cells.Item($r, 11).Value = 'Chat'
$ws.Cells.Item($r, 12).Value = 'N'
$ws.Cells.Item($r, 13).Value = 3
$ws.Cells.Item($r, 14).Value = 'OI'
$rowRange.Style = "Normal"
$ws.Rows.Item($r).EntireRow.AutoFit()

